$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '62.003.59'
Set-TextValue 'E2' '  -0.42%  '
Set-TextValue 'D3' '2.408.18'
Set-TextValue 'E3' '  -0.85%  '
Set-TextValue 'E4' '  -0.12%  '
Set-TextValue 'D5' '562.71'
Set-TextValue 'E5' '  +1.29%  '
Set-TextValue 'D6' '142.35'
Set-TextValue 'E6' '  -1.03%  '
Set-TextValue 'E7' '  +0.00%  '
Set-TextValue 'D8' '0.527'
Set-TextValue 'E8' '  -0.90%  '
Set-TextValue 'E9' '  +0.02%  '
Set-TextValue 'E10' '  -1.97%  '
Set-TextValue 'D11' '5.28'
Set-TextValue 'E11' '  -2.17%  '
Set-TextValue 'D12' '0.350'
Set-TextValue 'E12' '  -1.36%  '
Set-TextValue 'D13' '25.52'
Set-TextValue 'E13' '  -3.41%  '
Set-TextValue 'E14' '  -0.79%  '
Set-TextValue 'D15' '2.842.92'
Set-TextValue 'E15' '  -1.02%  '
Set-TextValue 'D16' '62.098.01'
Set-TextValue 'E16' '  +0.18%  '
Set-TextValue 'D17' '2.405.80'
Set-TextValue 'E17' '  -1.19%  '
Set-TextValue 'D18' '11.28'
Set-TextValue 'E18' '  +0.99%  '
Set-TextValue 'E19' '  +0.75%  '
Set-TextValue 'B20' 'BitcoinCash'
Set-TextValue 'C20' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D20' '320.97'
Set-TextValue 'E20' '  -1.17%  '
Set-TextValue 'B21' 'Polkadot'
Set-TextValue 'C21' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D21' '4.14'
Set-TextValue 'E21' '  -1.55%  '
Set-TextValue 'E22' '  -0.21%  '
Set-TextValue 'D23' '65.86'
Set-TextValue 'E23' '  +1.23%  '
Set-TextValue 'D24' '1.73'
Set-TextValue 'E24' '  -0.21%  '
Set-TextValue 'D25' '8.80'
Set-TextValue 'E25' '  -3.45%  '
Set-TextValue 'D26' '565.74'
Set-TextValue 'E26' '  +1.32%  '
Set-TextValue 'E27' '  +0.34%  '
Set-TextValue 'D28' '2.528.04'
Set-TextValue 'D29' '0.0₃0937'
Set-TextValue 'E29' '  -0.07%  '
Set-TextValue 'D30' '8.17'
Set-TextValue 'E30' '  -2.11%  '
Set-TextValue 'E31' '  -3.23%  '
Set-TextValue 'E32' '  -1.05%  '
Set-TextValue 'E33' '  +0.46%  '
Set-TextValue 'E34' '  -2.63%  '
Set-TextValue 'E35' '  -0.03%  '
Set-TextValue 'D36' '4.68'
Set-TextValue 'E36' '  -3.46%  '
Set-TextValue 'D37' '5.45'
Set-TextValue 'E37' '  -6.13%  '
Set-TextValue 'B38' 'Monero'
Set-TextValue 'C38' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D38' '151.90'
Set-TextValue 'E38' '  +3.78%  '
Set-TextValue 'B39' 'PolygonEcosystemToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D39' '0.379'
Set-TextValue 'E39' '  -1.58%  '
Set-TextValue 'D40' '18.60'
Set-TextValue 'E40' '  -0.96%  '
Set-TextValue 'D41' '1.78'
Set-TextValue 'E41' '  -9.83%  '
Set-TextValue 'E42' '  -0.03%  '
Set-TextValue 'E43' '  -0.58%  '
Set-TextValue 'D44' '147.57'
Set-TextValue 'E44' '  -2.06%  '
Set-TextValue 'D45' '3.62'
Set-TextValue 'E45' '  -0.64%  '
Set-TextValue 'E46' '  -2.14%  '
Set-TextValue 'D47' '19.83'
Set-TextValue 'E47' '  -2.79%  '
Set-TextValue 'D48' '0.590'
Set-TextValue 'E48' '  -0.35%  '
Set-TextValue 'E49' '  +0.41%  '
Set-TextValue 'E50' '  -1.26%  '
Set-TextValue 'E51' '  +0.14%  '
